$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'35.712.83"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.60%  '
$ws.Range('D3').Value = "'1.900.96"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.15%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = "'247.33"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.12%  '
$ws.Range('E6').Value = '  +0.24%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = "'43.27"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.58%  '
$ws.Range('D9').Value = "'57.50"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +10.63%  '
$ws.Range('E10').Value = '  +2.05%  '
$ws.Range('E11').Value = '  +2.25%  '
$ws.Range('D12').Value = "'0.0985"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.52%  '
$ws.Range('D13').Value = "'14.56"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +11.09%  '
$ws.Range('E14').Value = '  +11.66%  '
$ws.Range('D15').Value = "'2.173.66"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.04%  '
$ws.Range('E16').Value = '  +2.78%  '
$ws.Range('D17').Value = "'1.901.13"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.35%  '
$ws.Range('D18').Value = "'35.629.10"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.43%  '
$ws.Range('D19').Value = "'73.94"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.53%  '
$ws.Range('E20').Value = '  +1.44%  '
$ws.Range('D21').Value = "'247.63"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.11%  '
$ws.Range('D22').Value = "'13.05"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.71%  '
$ws.Range('E23').Value = '  +5.30%  '
$ws.Range('D24').Value = "'2.68"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +5.90%  '
$ws.Range('E26').Value = '  -1.20%  '
$ws.Range('D27').Value = "'166.94"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.75%  '
$ws.Range('E28').Value = '  +2.24%  '
$ws.Range('D29').Value = "'18.43"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.35%  '
$ws.Range('E30').Value = '  +0.80%  '
$ws.Range('D31').Value = "'4.40"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.57%  '
$ws.Range('D32').Value = "'0.0604"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.17%  '
$ws.Range('E33').Value = '  +0.73%  '
$ws.Range('B34').Value = 'BinanceUSD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D34').Value = "'1.00"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.05%  '
$ws.Range('B35').Value = 'WEMIXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').Value = "'1.84"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +13.02%  '
$ws.Range('E36').Value = '  -16.90%  '
$ws.Range('E37').Value = '  +0.38%  '
$ws.Range('D38').Value = "'1.98"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.08%  '
$ws.Range('D39').Value = "'0.0734"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +8.31%  '
$ws.Range('E40').Value = '  +6.66%  '
$ws.Range('D41').Value = "'99.73"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.56%  '
$ws.Range('D42').Value = "'17.19"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.30%  '
$ws.Range('B43').Value = 'Gas'
$ws.Range('C43').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D43').Value = "'14.52"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +20.03%  '
$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D44').Value = "'1.09"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.08%  '
$ws.Range('D45').Value = "'1.325.78"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.86%  '
$ws.Range('D46').Value = "'2.39"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.55%  '
$ws.Range('E47').Value = '  +1.27%  '
$ws.Range('E48').Value = '  -0.47%  '
$ws.Range('D49').Value = "'2.75"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.14%  '
$ws.Range('D50').Value = "'6.44"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.89%  '
$ws.Range('E51').Value = '  -0.94%  '
